$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.8621944130864
$ws.Range("B3").Value = 2.14537822100951
$ws.Range("B4").Value = 3.46552154686375
$ws.Range("B5").Value = 3.36913884643038
$ws.Range("B6").Value = 4.32730732243039
$ws.Range("B7").Value = 2.47780422113571
$ws.Range("B8").Value = 3.16733137281397
$ws.Range("B9").Value = 2.82216719124913
$ws.Range("B10").Value = 6.66388306797426
$ws.Range("B11").Value = 4.1536637478913
$ws.Range("B12").Value = 3.91530481763937
$ws.Range("B13").Value = 3.13792030753114
$ws.Range("B14").Value = 2.15276961793347
$ws.Range("B15").Value = 3.60897591672781
$ws.Range("B16").Value = 3.01202330711012
$ws.Range("B17").Value = 2.20115771408337
$ws.Range("B18").Value = 2.60225028486127
$ws.Range("B19").Value = 3.0954876204267
$ws.Range("B20").Value = 4.47119137181981
$ws.Range("B21").Value = 2.118383933431
$ws.Range("B22").Value = 3.64988170549473
$ws.Range("B23").Value = 2.76579225879593
$ws.Range("B24").Value = 3.17545574432584
$ws.Range("B25").Value = 2.17927627423665
$ws.Range("B26").Value = 4.34787976523288
$ws.Range("B27").Value = 3.00700347108456
$ws.Range("B28").Value = 2.02956718936119
$ws.Range("B29").Value = 2.31385640513292
$ws.Range("B30").Value = 3.77198592786264
$ws.Range("B31").Value = 1.88942284605715
$ws.Range("B32").Value = 3.65254531066039
$ws.Range("B33").Value = 4.15097670245584
$ws.Range("B34").Value = 4.09636956657109
$ws.Range("B35").Value = 3.61222594787496
$ws.Range("B36").Value = 2.28990438062146
$ws.Range("B37").Value = 3.20590948129276
$ws.Range("B38").Value = 3.30093086307198
$ws.Range("B39").Value = 2.59342364525698
$ws.Range("B40").Value = 2.87179293394553
$ws.Range("B41").Value = 3.33464271581237
$ws.Range("B42").Value = 3.45315157910072
$ws.Range("B43").Value = 2.2327698313403
$ws.Range("B44").Value = 3.22036310675211
$ws.Range("B45").Value = 4.29670849863915
$ws.Range("B46").Value = 2.1038507987104
$ws.Range("B47").Value = 1.86306300266008
$ws.Range("B48").Value = 2.9711141596086
$ws.Range("B49").Value = 2.3353454235078
$ws.Range("B50").Value = 2.79423210115508
$ws.Range("B51").Value = 2.3762366637113
$ws.Range("B52").Value = 1.89577685455957

Write-Output "Updated IER values for 51 states on sheet $($ws.Name)"
